$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Student full name: "Trần Nguyễn Kiên Tuấn" -> "Phạm Hoàng Anh"
# ------------------------------------------------------------------
$oldName = "Trần Nguyễn Kiên Tuấn"
$newName = "Phạm Hoàng Anh"
$full = $d.Content.Text
$idx = $full.IndexOf($oldName)
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $oldName.Length)
    $r.Text = $newName
}

# ------------------------------------------------------------------
# 2) Student ID: "BI12-468" -> "22BI13034" (keep the "Student ID: " label)
# ------------------------------------------------------------------
$oldId = "Student ID: BI12-468"
$newId = "Student ID: 22BI13034"
$full = $d.Content.Text
$idx = $full.IndexOf($oldId)
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $oldId.Length)
    $r.Text = $newId
}

# ------------------------------------------------------------------
# 3) Topic title: add "Movie and Cinema Management Application" right
#    after the existing "Topic: " label (same bold / size formatting).
#    NOTE: this paragraph lives inside a table cell, and collapsed
#    (zero-length) ranges misbehave for table content in this host,
#    so we replace the whole "Topic: " span (non-collapsed) instead
#    of inserting at a collapsed point.
# ------------------------------------------------------------------
$oldTopic = "Topic: "
$newTopic = "Topic: Movie and Cinema Management Application"
$full = $d.Content.Text
$idx = $full.IndexOf($oldTopic)
if ($idx -ge 0) {
    $r = $d.Range($idx, $idx + $oldTopic.Length)
    $r.Text = $newTopic
}
